$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2233009708737864
$ws.Range("C2").Value = 0.5194174757281553
$ws.Range("J2").Value = 0.004854368932038835
$ws.Range("P2").Value = 0.1407766990291262
$ws.Range("S2").Value = 0.1116504854368932
$ws.Range("C3").Value = 0.05454545454545454
$ws.Range("J3").Value = 0.02727272727272727
$ws.Range("P3").Value = 0.7909090909090909
$ws.Range("S3").Value = 0.1272727272727273
$ws.Range("P4").Value = 0.8095238095238095
$ws.Range("S4").Value = 0.1904761904761905
$ws.Range("B6").Value = 0.0797872340425532
$ws.Range("D6").Value = 0.01063829787234043
$ws.Range("F6").Value = 0.09042553191489362
$ws.Range("J6").Value = 0.2074468085106383
$ws.Range("O6").Value = 0.02659574468085106
$ws.Range("Q6").Value = 0.1595744680851064
$ws.Range("R6").Value = 0.03723404255319149
$ws.Range("S6").Value = 0.3882978723404255
$ws.Range("B7").Value = 0.0903954802259887
$ws.Range("D7").Value = 0.01694915254237288
$ws.Range("F7").Value = 0.04519774011299435
$ws.Range("J7").Value = 0.1468926553672316
$ws.Range("O7").Value = 0.005649717514124294
$ws.Range("Q7").Value = 0.192090395480226
$ws.Range("R7").Value = 0.07344632768361582
$ws.Range("S7").Value = 0.4293785310734463
$ws.Range("B8").Value = 0.05145413870246085
$ws.Range("D8").Value = 0.02237136465324385
$ws.Range("F8").Value = 0.058165548098434
$ws.Range("J8").Value = 0.1029082774049217
$ws.Range("O8").Value = 0.02013422818791946
$ws.Range("Q8").Value = 0.1677852348993289
$ws.Range("R8").Value = 0.1230425055928412
$ws.Range("S8").Value = 0.4541387024608501
$ws.Range("B9").Value = 0.075
$ws.Range("D9").Value = 0.04
$ws.Range("F9").Value = 0.04
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.025
$ws.Range("Q9").Value = 0.23
$ws.Range("R9").Value = 0.105
$ws.Range("S9").Value = 0.385
$ws.Range("B10").Value = 0.07964601769911504
$ws.Range("D10").Value = 0.01769911504424779
$ws.Range("F10").Value = 0.06371681415929203
$ws.Range("J10").Value = 0.1256637168141593
$ws.Range("O10").Value = 0.01858407079646018
$ws.Range("Q10").Value = 0.2353982300884956
$ws.Range("R10").Value = 0.08584070796460178
$ws.Range("S10").Value = 0.3734513274336283
$ws.Range("G11").Value = 0.1228813559322034
$ws.Range("J11").Value = 0.06779661016949153
$ws.Range("K11").Value = 0.1610169491525424
$ws.Range("L11").Value = 0.635593220338983
$ws.Range("S11").Value = 0.01271186440677966
$ws.Range("G12").Value = 0.7701863354037267
$ws.Range("J12").Value = 0.1490683229813665
$ws.Range("L12").Value = 0.06211180124223602
$ws.Range("S12").Value = 0.01863354037267081
$ws.Range("F13").Value = 0.025
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.225
$ws.Range("F15").Value = 0.0198019801980198
$ws.Range("H15").Value = 0.09900990099009901
$ws.Range("I15").Value = 0.0594059405940594
$ws.Range("J15").Value = 0.3762376237623762
$ws.Range("K15").Value = 0.07920792079207921
$ws.Range("M15").Value = 0.0198019801980198
$ws.Range("O15").Value = 0.09900990099009901
$ws.Range("S15").Value = 0.2475247524752475
$ws.Range("F16").Value = 0.01388888888888889
$ws.Range("H16").Value = 0.2291666666666667
$ws.Range("I16").Value = 0.09722222222222222
$ws.Range("J16").Value = 0.3888888888888889
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.02777777777777778
$ws.Range("O16").Value = 0.04861111111111111
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.01777777777777778
$ws.Range("H17").Value = 0.2111111111111111
$ws.Range("I17").Value = 0.1044444444444445
$ws.Range("J17").Value = 0.4222222222222222
$ws.Range("K17").Value = 0.08888888888888889
$ws.Range("M17").Value = 0.01555555555555556
$ws.Range("N17").Value = 0.002222222222222222
$ws.Range("O17").Value = 0.04666666666666667
$ws.Range("S17").Value = 0.09111111111111111
$ws.Range("F18").Value = 0.02631578947368421
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.1473684210526316
$ws.Range("J18").Value = 0.3789473684210526
$ws.Range("K18").Value = 0.08421052631578947
$ws.Range("M18").Value = 0.03157894736842105
$ws.Range("O18").Value = 0.04736842105263158
$ws.Range("S18").Value = 0.08421052631578947
$ws.Range("F19").Value = 0.01468048359240069
$ws.Range("H19").Value = 0.2288428324697755
$ws.Range("I19").Value = 0.08376511226252159
$ws.Range("J19").Value = 0.3652849740932643
$ws.Range("K19").Value = 0.09326424870466321
$ws.Range("M19").Value = 0.01727115716753022
$ws.Range("N19").Value = 0.0008635578583765112
$ws.Range("O19").Value = 0.0690846286701209
$ws.Range("S19").Value = 0.1269430051813472
